$wb = $excel.ActiveWorkbook

function Set-Cells {
    param($ws, $row, $values)
    foreach ($kv in $values.GetEnumerator()) {
        $ws.Cells.Item($row, $kv.Key).Value = $kv.Value
    }
}

$ws = $wb.Worksheets.Item("ALC")
Set-Cells $ws 19 @{8 = 2775.75; 9 = 2900.5; 10 = 2651; 11 = 2900.5; 12 = 2651; 13 = -2725.5; 14 = -3001}
Set-Cells $ws 28 @{8 = 340.35294; 9 = 272.7; 10 = 437; 11 = 272.7; 12 = 437; 13 = 212.3; 14 = -1407}
Set-Cells $ws 33 @{8 = 640.56525; 9 = 221.65; 10 = 3433.3333; 11 = 221.65; 12 = 3433.3333; 13 = 7.349999999999994; 14 = -3891.3333}
Set-Cells $ws 129 @{8 = 605.4231; 9 = 511.9524; 11 = 1535.8572; 13 = 3464.1428}
Set-Cells $ws 132 @{8 = 120084.33; 9 = 2494.0322; 10 = 437066.88; 11 = 7482.096600000001; 12 = 1311200.64; 13 = -4952.096600000001; 14 = -1316260.64}
Set-Cells $ws 137 @{8 = 46949.375; 9 = 91952.91; 10 = 8869.462; 11 = 275858.73; 12 = 26608.386; 13 = -273308.73; 14 = -31708.386}
Set-Cells $ws 138 @{8 = 1763.33; 9 = 915.8378; 10 = 2261.0635; 11 = 2747.5134; 12 = 6783.190500000001; 13 = 2392.4866; 14 = -17063.1905}

$ws = $wb.Worksheets.Item("ARM")
Set-Cells $ws 32 @{8 = 1909.09; 9 = 1717.7273; 10 = 3312.4167; 11 = 1717.7273; 12 = 3312.4167; 13 = -1430.7273; 14 = -3886.4167}
Set-Cells $ws 74 @{8 = 18224.492; 9 = 25073.857; 10 = 1302.5294; 11 = 25073.857; 12 = 1302.5294; 13 = -24199.857; 14 = -3050.5294}
Set-Cells $ws 77 @{8 = 18224.492; 9 = 25073.857; 10 = 1302.5294; 11 = 125369.285; 12 = 6512.646999999999; 13 = -121001.285; 14 = -15248.647}

$ws = $wb.Worksheets.Item("BSM")
Set-Cells $ws 107 @{8 = 683; 9 = 683; 11 = 683; 13 = 1237}
Set-Cells $ws 134 @{8 = 19297.016; 9 = 1025.0769; 10 = 124868.22; 11 = 3075.2307; 12 = 374604.66; 13 = -540.2307000000001; 14 = -379674.66}

$ws = $wb.Worksheets.Item("CRP")
Set-Cells $ws 16 @{8 = 1031.4445; 9 = 942.5; 10 = 1102.6; 11 = 942.5; 12 = 1102.6; 13 = -655.5; 14 = -1676.6}
Set-Cells $ws 99 @{8 = 4754.5454; 9 = 4362.5; 10 = 5800; 11 = 4362.5; 12 = 5800; 13 = -2864.5; 14 = -8796}
Set-Cells $ws 107 @{8 = 490.45834; 9 = 270.58334; 10 = 710.3333; 11 = 270.58334; 12 = 710.3333; 13 = 1649.41666; 14 = -4550.3333}
Set-Cells $ws 113 @{8 = 1031.4445; 9 = 942.5; 10 = 1102.6; 11 = 942.5; 12 = 1102.6; 13 = 1227.5; 14 = -5442.6}
Set-Cells $ws 126 @{8 = 4754.5454; 9 = 4362.5; 10 = 5800; 11 = 13087.5; 12 = 17400; 13 = -10617.5; 14 = -22340}

$ws = $wb.Worksheets.Item("CUL")
Set-Cells $ws 5 @{8 = 4778.08; 9 = 849.86664; 10 = 10670.4; 11 = 2549.59992; 12 = 32011.2; 13 = -2437.59992; 14 = -32235.2}
Set-Cells $ws 135 @{8 = 4778.08; 9 = 849.86664; 10 = 10670.4; 11 = 7648.79976; 12 = 96033.59999999999; 13 = -5113.79976; 14 = -101103.6}

$ws = $wb.Worksheets.Item("LTW")
Set-Cells $ws 7 @{8 = 3018.6562; 9 = 1916.95; 10 = 4854.8335; 11 = 1916.95; 12 = 4854.8335; 13 = -1804.95; 14 = -5078.8335}
Set-Cells $ws 40 @{8 = 65537.75; 9 = 2116; 11 = 2116; 13 = -1980}
Set-Cells $ws 61 @{8 = 1954.1364; 9 = 1543.1875; 10 = 3050; 11 = 1543.1875; 12 = 3050; 13 = -1341.1875; 14 = -3454}
Set-Cells $ws 113 @{8 = 1954.1364; 9 = 1543.1875; 10 = 3050; 11 = 1543.1875; 12 = 3050; 13 = 626.8125; 14 = -7390}
Set-Cells $ws 126 @{8 = 3018.6562; 9 = 1916.95; 10 = 4854.8335; 11 = 5750.85; 12 = 14564.5005; 13 = -3280.85; 14 = -19504.5005}

$ws = $wb.Worksheets.Item("WVR")
Set-Cells $ws 113 @{8 = 295.47058; 9 = 301.91666; 11 = 905.7499799999999; 13 = 1264.25002}
Set-Cells $ws 126 @{8 = 1039.2142; 9 = 882.7778; 10 = 1320.8; 11 = 2648.3334; 12 = 3962.4; 13 = -178.3334; 14 = -8902.4}

Write-Host "Applied market data updates to ALC, ARM, BSM, CRP, CUL, LTW, WVR sheets"
